# Daily attendance processing - 2026-01-06 10:37:28
# Applies the session-analysis refresh:
#  1) Column G "Recorded By" cells that list "System, <email>" are
#     re-ordered to "<email>, System" (72 recorded sessions).
#  2) The workbook-level Missing/Pending session counters (L7/L8) shift
#     by the 6 sessions that flipped from "Pending" to "Not Recorded".
#  3) The six group rows whose only pending session is now overdue move
#     from the Missing=2/Pending=8 bucket to Missing=3/Pending=7
#     (columns P/Q) for groups B1-10, B1-11, B1-12, B1-7, B1-8, B1-9.
#  4) The six corresponding per-session rows (30, 51, 72, 193, 214, 235)
#     flip their Status from "Pending" to "Not Recorded" and pick up the
#     same "Not Recorded" row formatting (pink fill) already used by the
#     other not-yet-due sessions such as row 10.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Normalise the "Recorded By" ordering for every already-recorded
#    session that still reads "System, dnasr281@gmail.com".
$recordedByRows = @(2,3,4,5,6,7,16,17,22,23,37,38,43,44,58,59,64,65,79,80,85,86,87,88,89,90,99,100,105,106,107,108,109,110,119,120,125,126,127,128,129,130,139,140,145,146,147,148,149,150,159,160,165,166,167,168,169,170,179,180,185,186,200,201,206,207,221,222,227,228,242,243)

foreach ($r in $recordedByRows) {
    $ws.Cells.Item($r, 7).Value = "dnasr281@gmail.com, System"
}

# 2) Roll-up "Missing Sessions" / "Pending Sessions" counters.
$ws.Range("L7").Value = 48
$ws.Range("L8").Value = 78

# 3) Per-group Missing/Pending shift for the six groups whose next
#    session (06/01/2026) just became overdue.
$groupStatRows = @(16,17,18,24,25,26)
foreach ($r in $groupStatRows) {
    $ws.Cells.Item($r, 16).Value = 3
    $ws.Cells.Item($r, 17).Value = 7
}

# 4) Flip the six now-overdue sessions from "Pending" to "Not Recorded",
#    matching the look of the other not-yet-recorded rows (copy the
#    formatting from row 10, a reference "Not Recorded" row).
$pendingToNotRecordedRows = @(30,51,72,193,214,235)
foreach ($r in $pendingToNotRecordedRows) {
    $ws.Range("A10:I10").Copy()
    $ws.Range("A" + $r + ":I" + $r).PasteSpecial(-4122)
    $ws.Cells.Item($r, 9).Value = "Not Recorded"
}
